$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-07-19 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-07-20 Thursday", 2) | Out-Null
$d.Content.Find.Execute("58×70=", $true, $false, $false, $false, $false, $true, 1, $false, "78×15=", 2) | Out-Null
$d.Content.Find.Execute("23×66=", $true, $false, $false, $false, $false, $true, 1, $false, "66×92=", 2) | Out-Null
$d.Content.Find.Execute("45×74=", $true, $false, $false, $false, $false, $true, 1, $false, "39×10=", 2) | Out-Null
$d.Content.Find.Execute("92×64=", $true, $false, $false, $false, $false, $true, 1, $false, "66×89=", 2) | Out-Null
$d.Content.Find.Execute("77×37=", $true, $false, $false, $false, $false, $true, 1, $false, "92×97=", 2) | Out-Null
$d.Content.Find.Execute("15×62=", $true, $false, $false, $false, $false, $true, 1, $false, "58×34=", 2) | Out-Null
$d.Content.Find.Execute("12×35=", $true, $false, $false, $false, $false, $true, 1, $false, "62×43=", 2) | Out-Null
$d.Content.Find.Execute("19×98=", $true, $false, $false, $false, $false, $true, 1, $false, "79×54=", 2) | Out-Null
$d.Content.Find.Execute("26×74=", $true, $false, $false, $false, $false, $true, 1, $false, "80×22=", 2) | Out-Null
$d.Content.Find.Execute("32×58=", $true, $false, $false, $false, $false, $true, 1, $false, "99×39=", 2) | Out-Null
$d.Content.Find.Execute("79×26=", $true, $false, $false, $false, $false, $true, 1, $false, "99×14=", 2) | Out-Null
$d.Content.Find.Execute("36×48=", $true, $false, $false, $false, $false, $true, 1, $false, "15×98=", 2) | Out-Null
$d.Content.Find.Execute("71×60=", $true, $false, $false, $false, $false, $true, 1, $false, "45×91=", 2) | Out-Null
$d.Content.Find.Execute("52×67=", $true, $false, $false, $false, $false, $true, 1, $false, "36×45=", 2) | Out-Null
$d.Content.Find.Execute("35×92=", $true, $false, $false, $false, $false, $true, 1, $false, "14×20=", 2) | Out-Null
$d.Content.Find.Execute("36×63=", $true, $false, $false, $false, $false, $true, 1, $false, "64×89=", 2) | Out-Null
$d.Content.Find.Execute("26×69=", $true, $false, $false, $false, $false, $true, 1, $false, "15×82=", 2) | Out-Null
$d.Content.Find.Execute("11×24=", $true, $false, $false, $false, $false, $true, 1, $false, "76×21=", 2) | Out-Null
$d.Content.Find.Execute("14×37=", $true, $false, $false, $false, $false, $true, 1, $false, "76×88=", 2) | Out-Null
$d.Content.Find.Execute("52×63=", $true, $false, $false, $false, $false, $true, 1, $false, "19×82=", 2) | Out-Null
$d.Content.Find.Execute("73×16=", $true, $false, $false, $false, $false, $true, 1, $false, "50×93=", 2) | Out-Null
$d.Content.Find.Execute("66×71=", $true, $false, $false, $false, $false, $true, 1, $false, "18×64=", 2) | Out-Null
$d.Content.Find.Execute("72×76=", $true, $false, $false, $false, $false, $true, 1, $false, "85×17=", 2) | Out-Null
$d.Content.Find.Execute("95×70=", $true, $false, $false, $false, $false, $true, 1, $false, "26×83=", 2) | Out-Null
$d.Content.Find.Execute("63×51=", $true, $false, $false, $false, $false, $true, 1, $false, "31×82=", 2) | Out-Null
$d.Content.Find.Execute("46×57=", $true, $false, $false, $false, $false, $true, 1, $false, "20×75=", 2) | Out-Null
$d.Content.Find.Execute("14×80=", $true, $false, $false, $false, $false, $true, 1, $false, "49×38=", 2) | Out-Null
$d.Content.Find.Execute("65×28=", $true, $false, $false, $false, $false, $true, 1, $false, "44×14=", 2) | Out-Null
$d.Content.Find.Execute("66×91=", $true, $false, $false, $false, $false, $true, 1, $false, "20×91=", 2) | Out-Null
$d.Content.Find.Execute("92×50=", $true, $false, $false, $false, $false, $true, 1, $false, "45×42=", 2) | Out-Null
$d.Content.Find.Execute("43×46=", $true, $false, $false, $false, $false, $true, 1, $false, "55×63=", 2) | Out-Null
$d.Content.Find.Execute("18×68=", $true, $false, $false, $false, $false, $true, 1, $false, "11×18=", 2) | Out-Null
$d.Content.Find.Execute("34×93=", $true, $false, $false, $false, $false, $true, 1, $false, "85×91=", 2) | Out-Null
$d.Content.Find.Execute("99×41=", $true, $false, $false, $false, $false, $true, 1, $false, "14×33=", 2) | Out-Null
$d.Content.Find.Execute("65×36=", $true, $false, $false, $false, $false, $true, 1, $false, "83×87=", 2) | Out-Null
$d.Content.Find.Execute("30×51=", $true, $false, $false, $false, $false, $true, 1, $false, "21×42=", 2) | Out-Null
$d.Content.Find.Execute("99×84=", $true, $false, $false, $false, $false, $true, 1, $false, "41×21=", 2) | Out-Null
$d.Content.Find.Execute("32×91=", $true, $false, $false, $false, $false, $true, 1, $false, "24×32=", 2) | Out-Null
$d.Content.Find.Execute("59×17=", $true, $false, $false, $false, $false, $true, 1, $false, "30×57=", 2) | Out-Null
$d.Content.Find.Execute("94×84=", $true, $false, $false, $false, $false, $true, 1, $false, "99×67=", 2) | Out-Null
$d.Content.Find.Execute("26×87=", $true, $false, $false, $false, $false, $true, 1, $false, "97×75=", 2) | Out-Null
$d.Content.Find.Execute("62×52=", $true, $false, $false, $false, $false, $true, 1, $false, "85×86=", 2) | Out-Null
$d.Content.Find.Execute("69×83=", $true, $false, $false, $false, $false, $true, 1, $false, "59×53=", 2) | Out-Null
$d.Content.Find.Execute("44×92=", $true, $false, $false, $false, $false, $true, 1, $false, "43×67=", 2) | Out-Null
$d.Content.Find.Execute("38×47=", $true, $false, $false, $false, $false, $true, 1, $false, "65×10=", 2) | Out-Null
$d.Content.Find.Execute("81×38=", $true, $false, $false, $false, $false, $true, 1, $false, "63×47=", 2) | Out-Null
$d.Content.Find.Execute("71×36=", $true, $false, $false, $false, $false, $true, 1, $false, "65×41=", 2) | Out-Null
$d.Content.Find.Execute("80×34=", $true, $false, $false, $false, $false, $true, 1, $false, "40×83=", 2) | Out-Null
$d.Content.Find.Execute("59×12=", $true, $false, $false, $false, $false, $true, 1, $false, "80×26=", 2) | Out-Null
$d.Content.Find.Execute("23×98=", $true, $false, $false, $false, $false, $true, 1, $false, "17×70=", 2) | Out-Null
$d.Content.Find.Execute("53×34=", $true, $false, $false, $false, $false, $true, 1, $false, "23×44=", 2) | Out-Null
$d.Content.Find.Execute("94×67=", $true, $false, $false, $false, $false, $true, 1, $false, "49×81=", 2) | Out-Null
$d.Content.Find.Execute("68×90=", $true, $false, $false, $false, $false, $true, 1, $false, "87×37=", 2) | Out-Null
$d.Content.Find.Execute("47×35=", $true, $false, $false, $false, $false, $true, 1, $false, "70×93=", 2) | Out-Null
$d.Content.Find.Execute("23×41=", $true, $false, $false, $false, $false, $true, 1, $false, "66×14=", 2) | Out-Null
$d.Content.Find.Execute("98×39=", $true, $false, $false, $false, $false, $true, 1, $false, "92×75=", 2) | Out-Null
$d.Content.Find.Execute("52×15=", $true, $false, $false, $false, $false, $true, 1, $false, "52×49=", 2) | Out-Null
$d.Content.Find.Execute("92×96=", $true, $false, $false, $false, $false, $true, 1, $false, "24×10=", 2) | Out-Null
$d.Content.Find.Execute("57×99=", $true, $false, $false, $false, $false, $true, 1, $false, "83×71=", 2) | Out-Null
$d.Content.Find.Execute("84×58=", $true, $false, $false, $false, $false, $true, 1, $false, "92×52=", 2) | Out-Null
$d.Content.Find.Execute("95×52=", $true, $false, $false, $false, $false, $true, 1, $false, "87×75=", 2) | Out-Null
$d.Content.Find.Execute("92×81=", $true, $false, $false, $false, $false, $true, 1, $false, "72×43=", 2) | Out-Null
$d.Content.Find.Execute("30×30=", $true, $false, $false, $false, $false, $true, 1, $false, "84×81=", 2) | Out-Null
$d.Content.Find.Execute("54×57=", $true, $false, $false, $false, $false, $true, 1, $false, "61×79=", 2) | Out-Null
$d.Content.Find.Execute("73×43=", $true, $false, $false, $false, $false, $true, 1, $false, "89×32=", 2) | Out-Null
$d.Content.Find.Execute("34×18=", $true, $false, $false, $false, $false, $true, 1, $false, "56×64=", 2) | Out-Null
$d.Content.Find.Execute("33×66=", $true, $false, $false, $false, $false, $true, 1, $false, "97×77=", 2) | Out-Null
$d.Content.Find.Execute("55×20=", $true, $false, $false, $false, $false, $true, 1, $false, "31×25=", 2) | Out-Null
$d.Content.Find.Execute("80×29=", $true, $false, $false, $false, $false, $true, 1, $false, "11×20=", 2) | Out-Null
$d.Content.Find.Execute("37×55=", $true, $false, $false, $false, $false, $true, 1, $false, "73×67=", 2) | Out-Null
$d.Content.Find.Execute("45×28=", $true, $false, $false, $false, $false, $true, 1, $false, "60×13=", 2) | Out-Null
$d.Content.Find.Execute("33×79=", $true, $false, $false, $false, $false, $true, 1, $false, "11×88=", 2) | Out-Null
$d.Content.Find.Execute("15×91=", $true, $false, $false, $false, $false, $true, 1, $false, "24×90=", 2) | Out-Null
$d.Content.Find.Execute("39×91=", $true, $false, $false, $false, $false, $true, 1, $false, "95×62=", 2) | Out-Null
$d.Content.Find.Execute("49×24=", $true, $false, $false, $false, $false, $true, 1, $false, "53×87=", 2) | Out-Null
$d.Content.Find.Execute("93×17=", $true, $false, $false, $false, $false, $true, 1, $false, "78×20=", 2) | Out-Null
$d.Content.Find.Execute("90×21=", $true, $false, $false, $false, $false, $true, 1, $false, "100×50=", 2) | Out-Null
$d.Content.Find.Execute("43×34=", $true, $false, $false, $false, $false, $true, 1, $false, "26×18=", 2) | Out-Null
$d.Content.Find.Execute("13×18=", $true, $false, $false, $false, $false, $true, 1, $false, "15×12=", 2) | Out-Null
$d.Content.Find.Execute("46×76=", $true, $false, $false, $false, $false, $true, 1, $false, "23×86=", 2) | Out-Null
$d.Content.Find.Execute("43×43=", $true, $false, $false, $false, $false, $true, 1, $false, "70×56=", 2) | Out-Null
$d.Content.Find.Execute("30×43=", $true, $false, $false, $false, $false, $true, 1, $false, "90×67=", 2) | Out-Null
$d.Content.Find.Execute("30×23=", $true, $false, $false, $false, $false, $true, 1, $false, "26×75=", 2) | Out-Null
$d.Content.Find.Execute("36×29=", $true, $false, $false, $false, $false, $true, 1, $false, "67×40=", 2) | Out-Null
$d.Content.Find.Execute("79×22=", $true, $false, $false, $false, $false, $true, 1, $false, "29×74=", 2) | Out-Null
$d.Content.Find.Execute("58×72=", $true, $false, $false, $false, $false, $true, 1, $false, "45×72=", 2) | Out-Null
$d.Content.Find.Execute("48×32=", $true, $false, $false, $false, $false, $true, 1, $false, "78×29=", 2) | Out-Null
$d.Content.Find.Execute("19×67=", $true, $false, $false, $false, $false, $true, 1, $false, "36×23=", 2) | Out-Null
$d.Content.Find.Execute("37×56=", $true, $false, $false, $false, $false, $true, 1, $false, "17×83=", 2) | Out-Null
$d.Content.Find.Execute("96×94=", $true, $false, $false, $false, $false, $true, 1, $false, "58×28=", 2) | Out-Null
$d.Content.Find.Execute("23×78=", $true, $false, $false, $false, $false, $true, 1, $false, "80×73=", 2) | Out-Null
$d.Content.Find.Execute("76×12=", $true, $false, $false, $false, $false, $true, 1, $false, "40×37=", 2) | Out-Null
$d.Content.Find.Execute("59×64=", $true, $false, $false, $false, $false, $true, 1, $false, "82×96=", 2) | Out-Null
$d.Content.Find.Execute("81×62=", $true, $false, $false, $false, $false, $true, 1, $false, "74×19=", 2) | Out-Null
$d.Content.Find.Execute("85×41=", $true, $false, $false, $false, $false, $true, 1, $false, "100×63=", 2) | Out-Null
$d.Content.Find.Execute("82×74=", $true, $false, $false, $false, $false, $true, 1, $false, "65×75=", 2) | Out-Null
$d.Content.Find.Execute("77×91=", $true, $false, $false, $false, $false, $true, 1, $false, "69×22=", 2) | Out-Null
$d.Content.Find.Execute("18×17=", $true, $false, $false, $false, $false, $true, 1, $false, "24×24=", 2) | Out-Null
$d.Content.Find.Execute("85×81=", $true, $false, $false, $false, $false, $true, 1, $false, "15×38=", 2) | Out-Null
$d.Content.Find.Execute("12×23=", $true, $false, $false, $false, $false, $true, 1, $false, "77×50=", 2) | Out-Null
